$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "group" column from A to B by inserting a new column at A.
$ws.Columns.Item(1).Insert()

# New column A header (Insert() already copied B1's formatting into A1)
$ws.Range("A1").Value = "Sample"

# New column A sample identifiers
$sampleIds = @(
    "SRR11856091",
    "SRR11856092",
    "SRR11856093",
    "SRR11856094",
    "SRR11856095",
    "SRR11856096",
    "SRR11856097",
    "SRR11856098",
    "SRR11856099",
    "SRR11856100"
)

for ($i = 0; $i -lt $sampleIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $sampleIds[$i]
}
